# Applies the "Added rule filters in UI" edit:
#  - Removes the ITERACJA 1 / ITERACJA 2 planning sections (paragraphs
#    27-41, 1-indexed) that were replaced by filter-related work.
#  - The trailing empty paragraph that used to close that block loses its
#    "Akapitzlist" list-paragraph formatting and becomes the new home of
#    the hidden "_GoBack" bookmark (Word automatically relocates that
#    bookmark to the most recent edit point, removing it from its old
#    location in the "postaci rownosci" paragraph).

$d = $word.ActiveDocument

# Range covering the whole ITERACJA 1 / ITERACJA 2 block, up to (but not
# including) the paragraph mark of the last paragraph in that block -
# that final mark/paragraph is kept and re-purposed below.
$blockStart = $d.Paragraphs(27).Range.Start
$blockEnd = $d.Paragraphs(41).Range.End
$block = $d.Range($blockStart, $blockEnd)
$block.Delete()

# The paragraph left behind (now empty) had "Akapitzlist" style only;
# resetting it to "Normal" drops its now-pointless pPr entirely.
$leftover = $d.Paragraphs(27)
$leftover.Range.Style = "Normal"

# Re-adding the "_GoBack" bookmark moves it here from wherever it used
# to be (there can only be one), matching Word's normal behaviour.
$d.Bookmarks.Add("_GoBack", $leftover.Range)
